$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Edit 1: paragraph that starts with "LUIS WRITE HERE Mario Party"
#   - drop the "LUIS WRITE HERE " placeholder run (with its bold /
#     italic / highlight / underline / big-size formatting)
#   - keep "Mario Party" (now plain, sz22/szCs18) as its own run
#   - append a new run continuing the sentence, same sz22/szCs18
# ------------------------------------------------------------------

$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "LUIS WRITE HERE") {
        $targetPara = $p
        break
    }
}

if ($targetPara -ne $null) {
    $full = $targetPara.Range
    # Exclude the trailing paragraph mark so the paragraph's own
    # <w:pPr> (sz22/szCs18 mark formatting) and identity survive.
    $bodyRange = $d.Range($full.Start, $full.End - 1)

    $xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="18"/></w:rPr><w:t>Mario Party</w:t></w:r><w:r><w:rPr><w:sz w:val="22"/><w:szCs w:val="18"/></w:rPr><w:t xml:space="preserve"> also inspired the idea of butter melt through its min-games. In Mario Party, there are thousands of mini-games that can be played by a user(s) as to determine if they won or not. Butter Melt is supposed to be similar to one of those mini-games but is intended to be played on mobile devices only.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $bodyRange.InsertXML($xml1)
}

# ------------------------------------------------------------------
# Edit 2: paragraph containing "Fjsekl;m" - strip the gramStart /
#   gramEnd proofErr markers that bracket the run (leave the run
#   itself untouched).
# ------------------------------------------------------------------

$targetPara2 = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "Fjsekl") {
        $targetPara2 = $p
        break
    }
}

if ($targetPara2 -ne $null) {
    $full2 = $targetPara2.Range

    $xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>Fjsekl;m</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

    $full2.InsertXML($xml2)
}

Write-Output "done"
